$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 304 entirely (the post "「大事なのは何者だったかじゃない。何者になるかだ」")
# This shifts all subsequent rows up by one.
$ws.Rows.Item(304).Delete()
